$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Committed Projects")

# Insert a new column before column I (9th column) and shift AREA..CULV_DURATION_N right.
$ws.Columns.Item(9).Insert()

# Populate the newly inserted column I with the PROJECTSOURCE data.
# (Order matters for shared-string table ordering: SAP, then ProjectBuilder,
# then the PROJECTSOURCE header, matching the original authoring order.)
$ws.Range("I2").Value = "SAP"
$ws.Range("I3").Value = "ProjectBuilder"
$ws.Range("I1").Value = "PROJECTSOURCE"
